$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("15÷9=1, 6", $true, $false, $false, $false, $false, $true, 0, $false, "12÷6=2, 0", 1) | Out-Null
$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("34÷7=4, 6", $true, $false, $false, $false, $false, $true, 0, $false, "18÷5=3, 3", 1) | Out-Null
$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("41÷6=6, 5", $true, $false, $false, $false, $false, $true, 0, $false, "99÷3=33, 0", 1) | Out-Null
$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("65÷6=10, 5", $true, $false, $false, $false, $false, $true, 0, $false, "13÷7=1, 6", 1) | Out-Null
$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("17÷8=2, 1", $true, $false, $false, $false, $false, $true, 0, $false, "68÷8=8, 4", 1) | Out-Null
$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("16÷2=8, 0", $true, $false, $false, $false, $false, $true, 0, $false, "45÷4=11, 1", 1) | Out-Null
$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("52÷6=8, 4", $true, $false, $false, $false, $false, $true, 0, $false, "92÷6=15, 2", 1) | Out-Null
$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("52÷6=8, 4", $true, $false, $false, $false, $false, $true, 0, $false, "75÷6=12, 3", 1) | Out-Null
$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("63÷4=15, 3", $true, $false, $false, $false, $false, $true, 0, $false, "90÷2=45, 0", 1) | Out-Null
$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("64÷8=8, 0", $true, $false, $false, $false, $false, $true, 0, $false, "96÷8=12, 0", 1) | Out-Null
$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("44÷5=8, 4", $true, $false, $false, $false, $false, $true, 0, $false, "61÷3=20, 1", 1) | Out-Null
$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("12÷2=6, 0", $true, $false, $false, $false, $false, $true, 0, $false, "35÷4=8, 3", 1) | Out-Null
$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("59÷7=8, 3", $true, $false, $false, $false, $false, $true, 0, $false, "59÷8=7, 3", 1) | Out-Null
$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("94÷6=15, 4", $true, $false, $false, $false, $false, $true, 0, $false, "43÷4=10, 3", 1) | Out-Null
$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("70÷6=11, 4", $true, $false, $false, $false, $false, $true, 0, $false, "48÷3=16, 0", 1) | Out-Null
$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("23÷3=7, 2", $true, $false, $false, $false, $false, $true, 0, $false, "30÷9=3, 3", 1) | Out-Null
$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("46÷5=9, 1", $true, $false, $false, $false, $false, $true, 0, $false, "59÷6=9, 5", 1) | Out-Null
$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("80÷6=13, 2", $true, $false, $false, $false, $false, $true, 0, $false, "39÷2=19, 1", 1) | Out-Null
$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("29÷5=5, 4", $true, $false, $false, $false, $false, $true, 0, $false, "63÷6=10, 3", 1) | Out-Null
$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("46÷6=7, 4", $true, $false, $false, $false, $false, $true, 0, $false, "66÷7=9, 3", 1) | Out-Null
$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("44÷3=14, 2", $true, $false, $false, $false, $false, $true, 0, $false, "31÷3=10, 1", 1) | Out-Null
$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("33÷3=11, 0", $true, $false, $false, $false, $false, $true, 0, $false, "66÷9=7, 3", 1) | Out-Null
$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("64÷7=9, 1", $true, $false, $false, $false, $false, $true, 0, $false, "38÷6=6, 2", 1) | Out-Null
$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("96÷6=16, 0", $true, $false, $false, $false, $false, $true, 0, $false, "31÷9=3, 4", 1) | Out-Null
$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("58÷2=29, 0", $true, $false, $false, $false, $false, $true, 0, $false, "43÷6=7, 1", 1) | Out-Null
